$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.876.22"
$ws.Range("E2").Value = "  -2.52%  "
$ws.Range("D3").Value = "2.287.11"
$ws.Range("E3").Value = "  -3.71%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.69"
$ws.Range("E5").Value = "  -0.78%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.74"
$ws.Range("E6").Value = "  -4.92%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.626"
$ws.Range("E7").Value = "  -1.48%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  -3.48%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.14"
$ws.Range("E10").Value = "  -6.96%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0901"
$ws.Range("E11").Value = "  -3.31%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.25"
$ws.Range("E12").Value = "  -4.16%  "
$ws.Range("E13").Value = "  -0.50%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.961"
$ws.Range("E14").Value = "  -5.25%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.22"
$ws.Range("E15").Value = "  -6.22%  "
$ws.Range("D16").Value = "2.633.19"
$ws.Range("E16").Value = "  -3.69%  "
$ws.Range("D17").Value = "2.285.22"
$ws.Range("E17").Value = "  -4.28%  "
$ws.Range("D18").Value = "41.984.98"
$ws.Range("E18").Value = "  -2.21%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.47"
$ws.Range("E19").Value = "  -3.22%  "
$ws.Range("E20").Value = "  -1.82%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "3.63"
$ws.Range("E21").Value = "  -2.61%  "
$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "282.65"
$ws.Range("E22").Value = "  +9.55%  "
$ws.Range("B23").Value = "Litecoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "73.31"
$ws.Range("E23").Value = "  -3.93%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.01"
$ws.Range("E24").Value = "  +5.70%  "
$ws.Range("E25").Value = "  -3.87%  "
$ws.Range("E26").Value = "  +0.55%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.77"
$ws.Range("E27").Value = "  -6.51%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.33"
$ws.Range("E28").Value = "  +3.19%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.92"
$ws.Range("E29").Value = "  -0.68%  "
$ws.Range("E30").Value = "  -3.98%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "163.36"
$ws.Range("E31").Value = "  -5.09%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0873"
$ws.Range("E32").Value = "  -3.07%  "
$ws.Range("E33").Value = "  -4.33%  "
$ws.Range("E34").Value = "  -3.42%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.134"
$ws.Range("E35").Value = "  +1.84%  "
$ws.Range("E36").Value = "  -5.39%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.51"
$ws.Range("E37").Value = "  -3.85%  "
$ws.Range("E38").Value = "  +6.62%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0347"
$ws.Range("E39").Value = "  -5.24%  "
$ws.Range("E40").Value = "  -7.00%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "100.07"
$ws.Range("E41").Value = "  +10.71%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.45"
$ws.Range("E42").Value = "  -5.43%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "69.40"
$ws.Range("E43").Value = "  -3.66%  "
$ws.Range("E44").Value = "  +0.16%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.223"
$ws.Range("E45").Value = "  -7.98%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "114.28"
$ws.Range("E46").Value = "  +0.41%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "11.85"
$ws.Range("E47").Value = "  -4.86%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "76.47"
$ws.Range("E48").Value = "  -1.27%  "
$ws.Range("E49").Value = "  -3.67%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.26"
$ws.Range("E50").Value = "  -5.75%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0992"
$ws.Range("E51").Value = "  -2.39%  "
